$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Change 1: merge the split "Dataset " / "\x96 meal" runs in the
#     POSTPARTUM sub-table header back into a single run, which also
#     removes the _GoBack bookmark that had been sitting between them. ---
$mealCell = $t.Cell(22, 3)
$mealCell.Range.Find.Execute("Dataset – meal", $false, $false, $false, $false, $false, $true, 1, $false, "Dataset – meal", 2) | Out-Null

# --- Change 2: drop the six detail rows for the "BMI group" breakdown
#     in the POSTPARTUM table (rows 23-28), then stamp a fresh, empty
#     _GoBack bookmark at the very start of the now-adjacent
#     "Return to early pregnancy weight" row's first cell. ---
for ($i = 0; $i -lt 6; $i++) {
    $t.Rows.Item(23).Delete()
}

$returnCell = $t.Cell(23, 1)
$startPos = $returnCell.Range.Start
$bmRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
